# "Agrego complemento, en el commit anterior no se adjunto"
# The previous commit's attachment (the real supplier rows) was missing,
# so the placeholder/sample rows that had been entered under the header
# of "Hoja1" (Base de datos Proveedores) are removed here, leaving only
# the header row until the real data is attached.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Header occupies row 1; delete every data row below it (rows 2..last).
$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -gt 1) {
    $ws.Rows("2:$lastRow").Delete()
}

# Restore the bottom pane's active cell/selection to B3, matching the
# cursor position left behind after clearing the sample rows.
[void]$ws.Range("B3").Select()
